$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rushing")

# Week 13 logging finished: M.Sargent did not play / has no stats this
# week, so remove his row entirely from the Rushing sheet. The rows below
# shift up to fill the gap, and the now-unused "M.Sargent" shared string
# is dropped automatically on save.
$ws.Rows.Item(4).Delete()

# Move the active selection to where the user's cursor ended up after the
# edit (mirrors the workbook's saved UI state).
$ws.Range("E17").Select() | Out-Null
